$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts everything below down by one row),
# carrying the new question "Sind Sie mobil" as a new parent question.
$ws.Rows("2:2").Insert()

# New row 2 content: ID = 1, Frage = "Sind Sie mobil"
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "Sind Sie mobil"
# The inserted row picked up the column B formatting (numFmt style) in B2;
# remove it completely since the new row has no ParentID value/format.
$ws.Range("B2").Clear()

# Renumber the ID column for what used to be rows 2-10 (now rows 3-11):
# they were 1..9, now they become 2..10.
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10

# The first two questions (rows 3 and 4) are now children of the new
# "Sind Sie mobil" question, so set their ParentID to 1.
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1

# A handful of ID cells lose their explicit wrap/valign style after the
# insert/renumber operation performed by Excel; restore the default style
# on those specific cells to match.
$ws.Range("A4").Style = "Standard"
$ws.Range("A6").Style = "Standard"
$ws.Range("A8").Style = "Standard"
$ws.Range("A10").Style = "Standard"

# The worksheet's saved sort range/condition also shifts down by one row
# along with the data it covers - re-apply it over its new location.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A3:A8"))
$sortObj.SetRange($ws.Range("A3:B8"))
$sortObj.Header = 0
$sortObj.Apply()

$ws.Range("B5").Select()
